# Applies three small wording/typo fixes to hw1.docx:
#   1. "-With all the c" + "ountries..." -> merge into a single run
#      "-With all the countries..." and drop the stray _GoBack bookmark
#      that used to sit between the two runs.
#   2. "pleaged" -> "pledged" (in "...backers_count and pleaged ")
#      and the _GoBack bookmark is re-anchored right after the
#      newly-typed "d" (i.e. between "pled" and "ged ").
#   3. "Pleged" -> "Pledged" (in "...quartile method to Pleged dataset...")

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge "-With all the c" / "ountries..." runs, drop bookmark
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "-With all the countries, theater has highest success rate, then music and film & video respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-With all the countries, theater has highest success rate, then music and film & video respectively.",
    2)

# ---------------------------------------------------------------------
# Change 2: "pleaged" -> "pledged", move the _GoBack bookmark so it
# sits right after the corrected "d" (between "pled" and "ged ").
# This text sits at the very end of its paragraph, so no barrier is
# needed to protect a trailing run.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    " and pleaged ", $true, $false, $false, $false, $false, $true, 1, $false,
    " and pledged ", 2)

$t = $d.Content.Text
$idx = $t.IndexOf(" and pledged ")
$splitPle = $idx + " and ple".Length
$splitD = $splitPle + 1

$rSplit1 = $d.Range($splitPle, $splitPle)
$d.Bookmarks.Add("TempSplitA", $rSplit1)
$rSplit2 = $d.Range($splitD, $splitD)
$d.Bookmarks.Add("_GoBack", $rSplit2)
$d.Bookmarks.Item("TempSplitA").Delete()

# ---------------------------------------------------------------------
# Change 3: "Pleged" -> "Pledged". The run after " dataset then " ("create
# a box plot. ") is untouched by the source edit, so pin it behind a
# throwaway bookmark first -- otherwise the engine's run-coalescing
# would merge it into the edited run since both share identical rPr.
# ---------------------------------------------------------------------
$t0 = $d.Content.Text
$idx0 = $t0.IndexOf("dataset then create")
$barrierPos = $idx0 + "dataset then ".Length
$rBarrier = $d.Range($barrierPos, $barrierPos)
$d.Bookmarks.Add("Barrier", $rBarrier)

$null = $d.Content.Find.Execute(
    "to Pleged dataset then", $true, $false, $false, $false, $false, $true, 1, $false,
    "to Pledged dataset then", 2)

$t2 = $d.Content.Text
$idx2 = $t2.IndexOf("to Pledged dataset")
$splitPle2 = $idx2 + "to Ple".Length
$splitD2 = $splitPle2 + 1

$rSplit3 = $d.Range($splitPle2, $splitPle2)
$d.Bookmarks.Add("TempSplitB", $rSplit3)
$rSplit4 = $d.Range($splitD2, $splitD2)
$d.Bookmarks.Add("TempSplitC", $rSplit4)
$d.Bookmarks.Item("TempSplitB").Delete()
$d.Bookmarks.Item("TempSplitC").Delete()
$d.Bookmarks.Item("Barrier").Delete()

Write-Output "done"
